# Apply the "add tradicional agent" simulation refresh:
# each row shifts to the next timestamp, the prior Real_Close becomes the new
# Trad_Prediction, and AI_Prediction is refreshed with new simulated output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45847.60416666666, 231.0700073242188, 231.1900024414062, 227.2487700910526),
    @(45847.64583333334, 232.0500030517578, 231.0700073242188, 225.8158038976075),
    @(45847.6875, 232.4015045166016, 232.0500030517578, 226.4340517966289),
    @(45847.72916666666, 232.8600006103516, 232.4015045166016, 231.8924958777332),
    @(45847.77083333334, 232.5731964111328, 232.8600006103516, 235.5706384732566),
    @(45847.8125, 231.8800048828125, 232.5731964111328, 232.8436628164346),
    @(45848.5625, 229.20849609375, 231.8800048828125, 226.224773302303),
    @(45848.60416666666, 230.0700073242188, 229.20849609375, 224.2264871855836),
    @(45848.64583333334, 229.8300018310547, 230.0700073242188, 225.4344317487278),
    @(45848.6875, 230.0399932861328, 229.8300018310547, 230.1408826919372),
    @(45848.72916666666, 229.8023986816406, 230.0399932861328, 224.7807957283161),
    @(45848.77083333334, 229.8699951171875, 229.8023986816406, 230.0397846004314),
    @(45848.8125, 229.9100036621094, 229.8699951171875, 237.05507798265),
    @(45849.5625, 230.6909027099609, 229.9100036621094, 231.5934559709231),
    @(45849.60416666666, 230.3200073242188, 230.6909027099609, 235.0129916631614),
    @(45849.64583333334, 231.5500030517578, 230.3200073242188, 220.6230412154993),
    @(45849.6875, 231.5700073242188, 231.5500030517578, 232.2085522082882),
    @(45849.72916666666, 231.25, 231.5700073242188, 234.9184864317685),
    @(45849.77083333334, 230.3600006103516, 231.25, 229.653119195338),
    @(45849.8125, 230.3600006103516, 230.3600006103516, 225.2614781701531),
    @(45852.5625, 228.0350036621094, 230.3600006103516, 230.1404909811272),
    @(45852.60416666666, 228.6999969482422, 228.0350036621094, 228.8489302236119),
    @(45852.64583333334, 228.5800018310547, 228.6999969482422, 226.6834112690531),
    @(45852.6875, 229.5899963378906, 228.5800018310547, 238.1465537965189),
    @(45852.72916666666, 229.0200042724609, 229.5899963378906, 220.5626103089149),
    @(45852.77083333334, 229.1049957275391, 229.0200042724609, 233.2263851522979),
    @(45852.8125, 228.6699981689453, 229.1049957275391, 224.7608805299284),
    @(45853.5625, 235.75, 228.6699981689453, 234.2132155466955),
    @(45853.60416666666, 236.8399963378906, 235.75, 237.3155049254995),
    @(45853.64583333334, 236.7299957275391, 236.8399963378906, 243.4179708486092),
    @(45853.6875, 236.0200042724609, 236.7299957275391, 244.0833319222015),
    @(45853.72916666666, 237.2599945068359, 236.0200042724609, 236.9673725024724),
    @(45853.77083333334, 237.2799987792969, 237.2599945068359, 245.3206976687475),
    @(45853.8125, 236.8780059814453, 237.2799987792969, 240.915593144717)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowValues[0]
    $ws.Cells.Item($r, 2).Value = $rowValues[1]
    $ws.Cells.Item($r, 3).Value = $rowValues[2]
    $ws.Cells.Item($r, 4).Value = $rowValues[3]
}

# The table now only spans 34 data rows (2-35); drop the now-obsolete last row
$ws.Rows.Item(36).Delete()

Write-Output "done"
